$d = $word.ActiveDocument

# --- Phase 1: text edits --------------------------------------------------
# The paragraph currently reads:
#   (run)"... моего отца - " (bookmark "_GoBack") (run)"нет."
# Trim the trailing " - " off the first run and prepend " \u2013 " (en dash)
# to the text that follows the bookmark, so the paragraph reads:
#   (run)"... моего отца" (bookmark) (run)" \u2013 нет."
$bm = $d.Bookmarks("_GoBack")
$bmStart = $bm.Start

# The three characters immediately before the bookmark are " - "
# (space, hyphen-minus, space); remove them.
$dashRange = $d.Range($bmStart - 3, $bmStart)
$dashRange.Text = ""

# Re-fetch the bookmark (its position shifts after the delete above) and
# rewrite the 4 characters right after it ("нет.") so they are prefixed
# with " \u2013 " (space, en dash, space).
$bm = $d.Bookmarks("_GoBack")
$bmEnd = $bm.End
$tailRange = $d.Range($bmEnd, $bmEnd + 4)
$enDash = [string][char]0x2013
$tailRange.Text = " " + $enDash + " " + $tailRange.Text

# The new tail text " \u2013 нет." is 7 characters long and starts right
# where the bookmark currently sits, so it now ends at bmEnd + 7.
$targetPos = $bmEnd + 7

# --- Phase 2: relocate the bookmark ---------------------------------------
# In the target document "_GoBack" sits at the very end of the paragraph
# (right after " \u2013 нет.", immediately before the paragraph mark)
# instead of between the two runs. Double-check targetPos really is that
# spot; fall back to searching the Paragraphs collection if not.
$bm = $d.Bookmarks("_GoBack")
$paras = $d.Paragraphs
$paraEnd = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Start -le $bm.Start -and $bm.Start -lt $p.Range.End) {
        $paraEnd = $p.Range.End
        break
    }
}
if ($paraEnd -ne -1 -and ($paraEnd - 1) -ne $targetPos) {
    $targetPos = $paraEnd - 1
}

$bm.Delete()

# Work around an engine quirk: Bookmarks.Add on a collapsed range placed
# exactly at (paragraph end - 1), i.e. immediately before the paragraph
# mark, mis-creates the bookmark (it ends up anchored at the wrong spot).
# Insert a throwaway placeholder character there first, anchor the new
# bookmark just before it, then remove the placeholder again.
$ph = $d.Range($targetPos, $targetPos)
$ph.InsertAfter("X")

$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$phRange = $d.Range($targetPos, $targetPos + 1)
$phRange.Text = ""
